$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.519.02"
$ws.Range("E2").Value = "  -2.03%  "

$ws.Range("D3").Value = "2.392.88"
$ws.Range("E3").Value = "  -2.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("D9").Value = "2.394.58"
$ws.Range("E9").Value = "  -1.76%  "

$ws.Range("E10").Value = "  -4.32%  "

$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("E12").Value = "  -1.84%  "

$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "2.813.62"

$ws.Range("E16").Value = "  -2.99%  "

$ws.Range("D17").Value = "60.597.38"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +14.78%  "

$ws.Range("D19").Value = "2.393.63"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.17"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.99%  "

$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "557.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -11.80%  "

$ws.Range("D29").Value = "2.522.62"
$ws.Range("E29").Value = "  -2.12%  "

$ws.Range("D30").Value = "0.0₃0913"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  -4.55%  "

$ws.Range("E33").Value = "  -3.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.52"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.25"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.08"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.64"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("E44").Value = "  -2.36%  "

$ws.Range("D45").Value = "0.0₆0281"
$ws.Range("E45").Value = "  -1.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.20"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.49"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.587"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("E49").Value = "  -2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.98"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.65%  "

$ws.Range("E51").Value = "  +0.34%  "
